# Rename the three inline logo pictures (two Pearson logos in the
# "first page" / "default" footers, one BTec logo in the "first page"
# header) by swapping their image2.png <-> image1.png / image1.jpg <->
# image2.jpg display names, per the target diff.
#
# wp:docPr/@name (InlineShape.Name) is the only part of the picture's
# display name that Word's automation model exposes for writing here;
# it is updated for each of the three inline pictures below.

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Footers -------------------------------------------------------
# Footer index 1 == default footer == footer2.xml == wp:docPr id="2"
#   Pearson logo: image2.png -> image1.png
$footer1 = $sec.Footers.Item(1)
$pearsonA = $footer1.Range.InlineShapes.Item(1)
$pearsonA.Name = "image1.png"

# Footer index 2 == first-page footer == footer1.xml == wp:docPr id="3"
#   Pearson logo: image2.png -> image1.png
$footer2 = $sec.Footers.Item(2)
$pearsonB = $footer2.Range.InlineShapes.Item(1)
$pearsonB.Name = "image1.png"

# --- Header ----------------------------------------------------------
# Header index 2 == first-page header == header1.xml == wp:docPr id="1"
#   BTec logo: image1.jpg -> image2.jpg
$header2 = $sec.Headers.Item(2)
$btec = $header2.Range.InlineShapes.Item(1)
$btec.Name = "image2.jpg"

Write-Output "Renamed 3 inline picture(s)."
